$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cells = $ws.Cells

# Row 32
$cells.Item(32, 1).Value = 112176074
$cells.Item(32, 2).Value = 98508
$cells.Item(32, 3).Value = 'Ovaliderad'
$cells.Item(32, 4).Value = 'LC'
$cells.Item(32, 5).Value = 1365
$cells.Item(32, 6).Value = 'Lappranunkel'
$cells.Item(32, 7).Value = 'Coptidium lapponicum'
$cells.Item(32, 8).Value = '(L.) Tzvelev'
$cells.Item(32, 9).Value = '''1'
$cells.Item(32, 16).Value = 'Håssjön, Ång'
$cells.Item(32, 17).Value = 602642.3927213018
$cells.Item(32, 18).Value = 7030561.374359156
$cells.Item(32, 19).Value = 5
$cells.Item(32, 20).Value = 'Västernorrland'
$cells.Item(32, 21).Value = 'Sollefteå'
$cells.Item(32, 22).Value = 'Ångermanland'
$cells.Item(32, 23).Value = 'Resele'
$cells.Item(32, 25).Value = '''2023-06-27'
$cells.Item(32, 26).Value = '00:00'
$cells.Item(32, 27).Value = '''2023-06-27'
$cells.Item(32, 28).Value = '00:00'
$cells.Item(32, 30).Value = $false
$cells.Item(32, 31).Value = $false
$cells.Item(32, 33).Value = $false
$cells.Item(32, 46).Formula = '=""'
$cells.Item(32, 49).Value = 'Pekka Bader'
$cells.Item(32, 50).Value = 'Pekka Bader, Anna-Maria Eriksson'
$cells.Item(32, 51).Value = 'Naturvärdesinventering Y-län'

# Row 33
$cells.Item(33, 1).Value = 112176087
$cells.Item(33, 2).Value = 89590
$cells.Item(33, 3).Value = 'Ovaliderad'
$cells.Item(33, 4).Value = 'VU'
$cells.Item(33, 5).Value = 48
$cells.Item(33, 6).Value = 'Lappticka'
$cells.Item(33, 7).Value = 'Amylocystis lapponica'
$cells.Item(33, 8).Value = '(Romell) Singer'
$cells.Item(33, 9).Value = '''1'
$cells.Item(33, 16).Value = 'Håssjön, Ång'
$cells.Item(33, 17).Value = 602805.8737273614
$cells.Item(33, 18).Value = 7030688.648317279
$cells.Item(33, 19).Value = 5
$cells.Item(33, 20).Value = 'Västernorrland'
$cells.Item(33, 21).Value = 'Sollefteå'
$cells.Item(33, 22).Value = 'Ångermanland'
$cells.Item(33, 23).Value = 'Resele'
$cells.Item(33, 25).Value = '''2023-06-27'
$cells.Item(33, 26).Value = '00:00'
$cells.Item(33, 27).Value = '''2023-06-27'
$cells.Item(33, 28).Value = '00:00'
$cells.Item(33, 30).Value = $false
$cells.Item(33, 31).Value = $false
$cells.Item(33, 33).Value = $false
$cells.Item(33, 46).Formula = '=""'
$cells.Item(33, 49).Value = 'Pekka Bader'
$cells.Item(33, 50).Value = 'Pekka Bader, Anna-Maria Eriksson'
$cells.Item(33, 51).Value = 'Naturvärdesinventering Y-län'

# Row 34
$cells.Item(34, 1).Value = 112176102
$cells.Item(34, 2).Value = 89686
$cells.Item(34, 3).Value = 'Ovaliderad'
$cells.Item(34, 4).Value = 'NT'
$cells.Item(34, 5).Value = 658
$cells.Item(34, 6).Value = 'Rosenticka'
$cells.Item(34, 7).Value = 'Rhodofomes roseus'
$cells.Item(34, 8).Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$cells.Item(34, 9).Value = '''1'
$cells.Item(34, 16).Value = 'Håssjön, Ång'
$cells.Item(34, 17).Value = 602660.3503862405
$cells.Item(34, 18).Value = 7030715.82964331
$cells.Item(34, 19).Value = 5
$cells.Item(34, 20).Value = 'Västernorrland'
$cells.Item(34, 21).Value = 'Sollefteå'
$cells.Item(34, 22).Value = 'Ångermanland'
$cells.Item(34, 23).Value = 'Resele'
$cells.Item(34, 25).Value = '''2023-06-14'
$cells.Item(34, 26).Value = '00:00'
$cells.Item(34, 27).Value = '''2023-06-14'
$cells.Item(34, 28).Value = '00:00'
$cells.Item(34, 30).Value = $false
$cells.Item(34, 31).Value = $false
$cells.Item(34, 33).Value = $false
$cells.Item(34, 46).Formula = '=""'
$cells.Item(34, 49).Value = 'Pekka Bader'
$cells.Item(34, 50).Value = 'Pekka Bader'
$cells.Item(34, 51).Value = 'Naturvärdesinventering Y-län'

# Row 35
$cells.Item(35, 1).Value = 112176069
$cells.Item(35, 2).Value = 18534
$cells.Item(35, 3).Value = 'Ovaliderad'
$cells.Item(35, 4).Value = 'EN'
$cells.Item(35, 5).Value = 101797
$cells.Item(35, 6).Value = 'Karelsk barkfluga'
$cells.Item(35, 7).Value = 'Xylomya czekanovskii'
$cells.Item(35, 8).Value = 'Pleske, 1925'
$cells.Item(35, 9).Value = '''1'
$cells.Item(35, 11).Value = 'puppa'
$cells.Item(35, 16).Value = 'Håssjön, Ång'
$cells.Item(35, 17).Value = 602774.5095870381
$cells.Item(35, 18).Value = 7030644.124442326
$cells.Item(35, 19).Value = 5
$cells.Item(35, 20).Value = 'Västernorrland'
$cells.Item(35, 21).Value = 'Sollefteå'
$cells.Item(35, 22).Value = 'Ångermanland'
$cells.Item(35, 23).Value = 'Resele'
$cells.Item(35, 25).Value = '''2023-06-27'
$cells.Item(35, 26).Value = '00:00'
$cells.Item(35, 27).Value = '''2023-06-27'
$cells.Item(35, 28).Value = '00:00'
$cells.Item(35, 30).Value = $false
$cells.Item(35, 31).Value = $false
$cells.Item(35, 33).Value = $false
$cells.Item(35, 44).Formula = '=""'
$cells.Item(35, 46).Formula = '=""'
$cells.Item(35, 49).Value = 'Pekka Bader'
$cells.Item(35, 50).Value = 'Pekka Bader, Anna-Maria Eriksson'
$cells.Item(35, 51).Value = 'Naturvärdesinventering Y-län'

# Row 36
$cells.Item(36, 1).Value = 112176095
$cells.Item(36, 2).Value = 89845
$cells.Item(36, 3).Value = 'Ovaliderad'
$cells.Item(36, 4).Value = 'VU'
$cells.Item(36, 5).Value = 1209
$cells.Item(36, 6).Value = 'Rynkskinn'
$cells.Item(36, 7).Value = 'Phlebia centrifuga'
$cells.Item(36, 8).Value = 'P.Karst.'
$cells.Item(36, 9).Value = '''1'
$cells.Item(36, 16).Value = 'Håssjön, Ång'
$cells.Item(36, 17).Value = 602796.3254976775
$cells.Item(36, 18).Value = 7030566.314026224
$cells.Item(36, 19).Value = 5
$cells.Item(36, 20).Value = 'Västernorrland'
$cells.Item(36, 21).Value = 'Sollefteå'
$cells.Item(36, 22).Value = 'Ångermanland'
$cells.Item(36, 23).Value = 'Resele'
$cells.Item(36, 25).Value = '''2023-06-14'
$cells.Item(36, 26).Value = '00:00'
$cells.Item(36, 27).Value = '''2023-06-14'
$cells.Item(36, 28).Value = '00:00'
$cells.Item(36, 30).Value = $false
$cells.Item(36, 31).Value = $false
$cells.Item(36, 33).Value = $false
$cells.Item(36, 46).Formula = '=""'
$cells.Item(36, 49).Value = 'Pekka Bader'
$cells.Item(36, 50).Value = 'Pekka Bader'
$cells.Item(36, 51).Value = 'Naturvärdesinventering Y-län'

# Row 37
$cells.Item(37, 1).Value = 112176096
$cells.Item(37, 2).Value = 12450
$cells.Item(37, 3).Value = 'Ovaliderad'
$cells.Item(37, 4).Value = 'EN'
$cells.Item(37, 5).Value = 101692
$cells.Item(37, 6).Value = 'Större barkplattbagge'
$cells.Item(37, 7).Value = 'Pytho kolwensis'
$cells.Item(37, 8).Value = 'Sahlberg, 1833'
$cells.Item(37, 9).Value = '''1'
$cells.Item(37, 11).Value = 'larv/nymf'
$cells.Item(37, 16).Value = 'Håssjön, Ång'
$cells.Item(37, 17).Value = 602868.7858234661
$cells.Item(37, 18).Value = 7030590.176470381
$cells.Item(37, 19).Value = 5
$cells.Item(37, 20).Value = 'Västernorrland'
$cells.Item(37, 21).Value = 'Sollefteå'
$cells.Item(37, 22).Value = 'Ångermanland'
$cells.Item(37, 23).Value = 'Resele'
$cells.Item(37, 25).Value = '''2023-06-14'
$cells.Item(37, 26).Value = '00:00'
$cells.Item(37, 27).Value = '''2023-06-14'
$cells.Item(37, 28).Value = '00:00'
$cells.Item(37, 29).Value = 'larv 20-25 mm'
$cells.Item(37, 30).Value = $false
$cells.Item(37, 31).Value = $false
$cells.Item(37, 33).Value = $false
$cells.Item(37, 46).Formula = '=""'
$cells.Item(37, 49).Value = 'Pekka Bader'
$cells.Item(37, 50).Value = 'Pekka Bader'
$cells.Item(37, 51).Value = 'Naturvärdesinventering Y-län'

# Row 38
$cells.Item(38, 1).Value = 112176108
$cells.Item(38, 2).Value = 89686
$cells.Item(38, 3).Value = 'Ovaliderad'
$cells.Item(38, 4).Value = 'NT'
$cells.Item(38, 5).Value = 658
$cells.Item(38, 6).Value = 'Rosenticka'
$cells.Item(38, 7).Value = 'Rhodofomes roseus'
$cells.Item(38, 8).Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$cells.Item(38, 9).Value = '''1'
$cells.Item(38, 16).Value = 'Håssjön, Ång'
$cells.Item(38, 17).Value = 602830.8989684412
$cells.Item(38, 18).Value = 7030664.777842449
$cells.Item(38, 19).Value = 5
$cells.Item(38, 20).Value = 'Västernorrland'
$cells.Item(38, 21).Value = 'Sollefteå'
$cells.Item(38, 22).Value = 'Ångermanland'
$cells.Item(38, 23).Value = 'Resele'
$cells.Item(38, 25).Value = '''2023-06-14'
$cells.Item(38, 26).Value = '00:00'
$cells.Item(38, 27).Value = '''2023-06-14'
$cells.Item(38, 28).Value = '00:00'
$cells.Item(38, 30).Value = $false
$cells.Item(38, 31).Value = $false
$cells.Item(38, 33).Value = $false
$cells.Item(38, 46).Formula = '=""'
$cells.Item(38, 49).Value = 'Pekka Bader'
$cells.Item(38, 50).Value = 'Pekka Bader'
$cells.Item(38, 51).Value = 'Naturvärdesinventering Y-län'

# Row 39
$cells.Item(39, 1).Value = 112176088
$cells.Item(39, 2).Value = 89425
$cells.Item(39, 3).Value = 'Ovaliderad'
$cells.Item(39, 4).Value = 'NT'
$cells.Item(39, 5).Value = 5442
$cells.Item(39, 6).Value = 'Tallticka'
$cells.Item(39, 7).Value = 'Porodaedalea pini'
$cells.Item(39, 8).Value = '(Brot.) Murrill'
$cells.Item(39, 9).Value = '''1'
$cells.Item(39, 16).Value = 'Håssjön, Ång'
$cells.Item(39, 17).Value = 602858.8694568657
$cells.Item(39, 18).Value = 7030591.203716032
$cells.Item(39, 19).Value = 5
$cells.Item(39, 20).Value = 'Västernorrland'
$cells.Item(39, 21).Value = 'Sollefteå'
$cells.Item(39, 22).Value = 'Ångermanland'
$cells.Item(39, 23).Value = 'Resele'
$cells.Item(39, 25).Value = '''2023-06-27'
$cells.Item(39, 26).Value = '00:00'
$cells.Item(39, 27).Value = '''2023-06-27'
$cells.Item(39, 28).Value = '00:00'
$cells.Item(39, 30).Value = $false
$cells.Item(39, 31).Value = $false
$cells.Item(39, 33).Value = $false
$cells.Item(39, 46).Formula = '=""'
$cells.Item(39, 49).Value = 'Pekka Bader'
$cells.Item(39, 50).Value = 'Pekka Bader, Anna-Maria Eriksson'
$cells.Item(39, 51).Value = 'Naturvärdesinventering Y-län'

# Row 40
$cells.Item(40, 1).Value = 112176093
$cells.Item(40, 2).Value = 6202
$cells.Item(40, 3).Value = 'Ovaliderad'
$cells.Item(40, 4).Value = 'LC'
$cells.Item(40, 5).Value = 105336
$cells.Item(40, 6).Value = 'Vanlig flatbagge'
$cells.Item(40, 7).Value = 'Peltis ferruginea'
$cells.Item(40, 8).Value = '(Linnaeus, 1758)'
$cells.Item(40, 9).Value = '''1'
$cells.Item(40, 16).Value = 'Håssjön, Ång'
$cells.Item(40, 17).Value = 602865.1212630216
$cells.Item(40, 18).Value = 7030578.394125014
$cells.Item(40, 19).Value = 5
$cells.Item(40, 20).Value = 'Västernorrland'
$cells.Item(40, 21).Value = 'Sollefteå'
$cells.Item(40, 22).Value = 'Ångermanland'
$cells.Item(40, 23).Value = 'Resele'
$cells.Item(40, 25).Value = '''2023-06-22'
$cells.Item(40, 26).Value = '00:00'
$cells.Item(40, 27).Value = '''2023-06-22'
$cells.Item(40, 28).Value = '00:00'
$cells.Item(40, 30).Value = $false
$cells.Item(40, 31).Value = $false
$cells.Item(40, 33).Value = $false
$cells.Item(40, 46).Formula = '=""'
$cells.Item(40, 49).Value = 'Pekka Bader'
$cells.Item(40, 50).Value = 'Pekka Bader'
$cells.Item(40, 51).Value = 'Naturvärdesinventering Y-län'

